$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3414.5715
$ws.Cells.Item(64, 9).Value = 3345.8635
$ws.Cells.Item(64, 10).Value = 3666.5
$ws.Cells.Item(64, 11).Value = 3345.8635
$ws.Cells.Item(64, 12).Value = 3666.5
$ws.Cells.Item(64, 13).Value = -3097.8635
$ws.Cells.Item(64, 14).Value = -4162.5
$ws.Cells.Item(67, 8).Value = 3414.5715
$ws.Cells.Item(67, 9).Value = 3345.8635
$ws.Cells.Item(67, 10).Value = 3666.5
$ws.Cells.Item(67, 11).Value = 3345.8635
$ws.Cells.Item(67, 12).Value = 3666.5
$ws.Cells.Item(67, 13).Value = -2487.8635
$ws.Cells.Item(67, 14).Value = -5382.5
$ws.Cells.Item(76, 8).Value = 3363.9832
$ws.Cells.Item(76, 9).Value = 3275.745
$ws.Cells.Item(76, 10).Value = 3926.5
$ws.Cells.Item(76, 11).Value = 3275.745
$ws.Cells.Item(76, 12).Value = 3926.5
$ws.Cells.Item(76, 13).Value = -2960.745
$ws.Cells.Item(76, 14).Value = -4556.5
$ws.Cells.Item(79, 8).Value = 3363.9832
$ws.Cells.Item(79, 9).Value = 3275.745
$ws.Cells.Item(79, 10).Value = 3926.5
$ws.Cells.Item(79, 11).Value = 3275.745
$ws.Cells.Item(79, 12).Value = 3926.5
$ws.Cells.Item(79, 13).Value = -2183.745
$ws.Cells.Item(79, 14).Value = -6110.5
$ws.Cells.Item(86, 8).Value = 33176.188
$ws.Cells.Item(86, 9).Value = 78591.234
$ws.Cells.Item(86, 10).Value = 2102.7368
$ws.Cells.Item(86, 11).Value = 78591.234
$ws.Cells.Item(86, 12).Value = 2102.7368
$ws.Cells.Item(86, 13).Value = -77468.234
$ws.Cells.Item(86, 14).Value = -4348.736800000001
$ws.Cells.Item(89, 8).Value = 33176.188
$ws.Cells.Item(89, 9).Value = 78591.234
$ws.Cells.Item(89, 10).Value = 2102.7368
$ws.Cells.Item(89, 11).Value = 392956.17
$ws.Cells.Item(89, 12).Value = 10513.684
$ws.Cells.Item(89, 13).Value = -387340.17
$ws.Cells.Item(89, 14).Value = -21745.684
$ws.Cells.Item(106, 8).Value = 2498.182
$ws.Cells.Item(106, 9).Value = 1493.3334
$ws.Cells.Item(106, 10).Value = 2875
$ws.Cells.Item(106, 11).Value = 1493.3334
$ws.Cells.Item(106, 12).Value = 2875
$ws.Cells.Item(106, 13).Value = -862.3334
$ws.Cells.Item(106, 14).Value = -4137
$ws.Cells.Item(137, 8).Value = 2454.541
$ws.Cells.Item(137, 9).Value = 1487.8125
$ws.Cells.Item(137, 11).Value = 4463.4375
$ws.Cells.Item(137, 13).Value = -1913.4375
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4140.18
$ws.Cells.Item(32, 9).Value = 2936.7952
$ws.Cells.Item(32, 10).Value = 10015.529
$ws.Cells.Item(32, 11).Value = 2936.7952
$ws.Cells.Item(32, 12).Value = 10015.529
$ws.Cells.Item(32, 13).Value = -2649.7952
$ws.Cells.Item(32, 14).Value = -10589.529
$ws.Cells.Item(63, 8).Value = 2002.5
$ws.Cells.Item(63, 9).Value = 2002.5
$ws.Cells.Item(63, 11).Value = 2002.5
$ws.Cells.Item(63, 13).Value = -1316.5
$ws.Cells.Item(66, 8).Value = 2002.5
$ws.Cells.Item(66, 9).Value = 2002.5
$ws.Cells.Item(66, 11).Value = 10012.5
$ws.Cells.Item(66, 13).Value = -6580.5
$ws.Cells.Item(74, 8).Value = 5990.839
$ws.Cells.Item(74, 9).Value = 3382.2856
$ws.Cells.Item(74, 11).Value = 3382.2856
$ws.Cells.Item(74, 13).Value = -2508.2856
$ws.Cells.Item(77, 8).Value = 5990.839
$ws.Cells.Item(77, 9).Value = 3382.2856
$ws.Cells.Item(77, 11).Value = 16911.428
$ws.Cells.Item(77, 13).Value = -12543.428
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 8773977
$ws.Cells.Item(86, 9).Value = 10103128
$ws.Cells.Item(86, 10).Value = 1581.2
$ws.Cells.Item(86, 11).Value = 10103128
$ws.Cells.Item(86, 12).Value = 1581.2
$ws.Cells.Item(86, 13).Value = -10102005
$ws.Cells.Item(86, 14).Value = -3827.2
$ws.Cells.Item(89, 8).Value = 8773977
$ws.Cells.Item(89, 9).Value = 10103128
$ws.Cells.Item(89, 10).Value = 1581.2
$ws.Cells.Item(89, 11).Value = 50515640
$ws.Cells.Item(89, 12).Value = 7906
$ws.Cells.Item(89, 13).Value = -50510024
$ws.Cells.Item(89, 14).Value = -19138
$ws.Cells.Item(105, 8).Value = 4698.952
$ws.Cells.Item(105, 9).Value = 4168.0386
$ws.Cells.Item(105, 10).Value = 5561.6875
$ws.Cells.Item(105, 11).Value = 4168.0386
$ws.Cells.Item(105, 12).Value = 5561.6875
$ws.Cells.Item(105, 13).Value = -2421.0386
$ws.Cells.Item(105, 14).Value = -9055.6875
$ws.Cells.Item(134, 8).Value = 5797.731
$ws.Cells.Item(134, 9).Value = 5267
$ws.Cells.Item(134, 10).Value = 9866.666999999999
$ws.Cells.Item(134, 11).Value = 15801
$ws.Cells.Item(134, 12).Value = 29600.001
$ws.Cells.Item(134, 13).Value = -13266
$ws.Cells.Item(134, 14).Value = -34670.001
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3519.375
$ws.Cells.Item(62, 9).Value = 2933.3333
$ws.Cells.Item(62, 11).Value = 2933.3333
$ws.Cells.Item(62, 13).Value = -2309.3333
$ws.Cells.Item(65, 8).Value = 3519.375
$ws.Cells.Item(65, 9).Value = 2933.3333
$ws.Cells.Item(65, 11).Value = 14666.6665
$ws.Cells.Item(65, 13).Value = -11546.6665
$ws.Cells.Item(109, 8).Value = 35000
$ws.Cells.Item(109, 10).Value = 35000
$ws.Cells.Item(109, 12).Value = 35000
$ws.Cells.Item(109, 14).Value = -37080
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 25641236
$ws.Cells.Item(12, 9).Value = 76923270
$ws.Cells.Item(12, 10).Value = 217.23077
$ws.Cells.Item(12, 11).Value = 230769810
$ws.Cells.Item(12, 12).Value = 651.69231
$ws.Cells.Item(12, 13).Value = -230769637
$ws.Cells.Item(12, 14).Value = -997.69231
$ws.Cells.Item(22, 8).Value = 1310
$ws.Cells.Item(22, 9).Value = 1025
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 11).Value = 3075
$ws.Cells.Item(22, 12).Value = 4500
$ws.Cells.Item(22, 13).Value = -2906
$ws.Cells.Item(22, 14).Value = -4838
$ws.Cells.Item(27, 8).Value = 1310
$ws.Cells.Item(27, 9).Value = 1025
$ws.Cells.Item(27, 10).Value = 1500
$ws.Cells.Item(27, 11).Value = 3075
$ws.Cells.Item(27, 12).Value = 4500
$ws.Cells.Item(27, 13).Value = -2973
$ws.Cells.Item(27, 14).Value = -4704
$ws.Cells.Item(34, 8).Value = 2440
$ws.Cells.Item(34, 9).Value = 229
$ws.Cells.Item(34, 10).Value = 3740.5881
$ws.Cells.Item(34, 11).Value = 687
$ws.Cells.Item(34, 12).Value = 11221.7643
$ws.Cells.Item(34, 13).Value = -603
$ws.Cells.Item(34, 14).Value = -11389.7643
$ws.Cells.Item(40, 8).Value = 87.36364
$ws.Cells.Item(40, 9).Value = 87.59999999999999
$ws.Cells.Item(40, 10).Value = 85
$ws.Cells.Item(40, 11).Value = 350.4
$ws.Cells.Item(40, 12).Value = 340
$ws.Cells.Item(40, 13).Value = -281.4
$ws.Cells.Item(40, 14).Value = -478
$ws.Cells.Item(46, 8).Value = 2328.5715
$ws.Cells.Item(46, 9).Value = 200
$ws.Cells.Item(46, 10).Value = 3180
$ws.Cells.Item(46, 11).Value = 600
$ws.Cells.Item(46, 12).Value = 9540
$ws.Cells.Item(46, 13).Value = -509
$ws.Cells.Item(46, 14).Value = -9722
$ws.Cells.Item(68, 8).Value = 6524.5884
$ws.Cells.Item(68, 9).Value = 580.2
$ws.Cells.Item(68, 10).Value = 9001.416999999999
$ws.Cells.Item(68, 11).Value = 1740.6
$ws.Cells.Item(68, 12).Value = 27004.251
$ws.Cells.Item(68, 13).Value = -929.6000000000001
$ws.Cells.Item(68, 14).Value = -28626.251
$ws.Cells.Item(71, 8).Value = 6524.5884
$ws.Cells.Item(71, 9).Value = 580.2
$ws.Cells.Item(71, 10).Value = 9001.416999999999
$ws.Cells.Item(71, 11).Value = 5221.8
$ws.Cells.Item(71, 12).Value = 81012.753
$ws.Cells.Item(71, 13).Value = -1165.8
$ws.Cells.Item(71, 14).Value = -89124.753
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5716.9785
$ws.Cells.Item(70, 9).Value = 5676.579
$ws.Cells.Item(70, 10).Value = 5744.393
$ws.Cells.Item(70, 11).Value = 5676.579
$ws.Cells.Item(70, 12).Value = 5744.393
$ws.Cells.Item(70, 13).Value = -5406.579
$ws.Cells.Item(70, 14).Value = -6284.393
$ws.Cells.Item(73, 8).Value = 5716.9785
$ws.Cells.Item(73, 9).Value = 5676.579
$ws.Cells.Item(73, 10).Value = 5744.393
$ws.Cells.Item(73, 11).Value = 5676.579
$ws.Cells.Item(73, 12).Value = 5744.393
$ws.Cells.Item(73, 13).Value = -4740.579
$ws.Cells.Item(73, 14).Value = -7616.393
$ws.Cells.Item(80, 8).Value = 8550
$ws.Cells.Item(80, 10).Value = 4000
$ws.Cells.Item(80, 12).Value = 4000
$ws.Cells.Item(80, 14).Value = -5996
$ws.Cells.Item(83, 8).Value = 8550
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 12).Value = 20000
$ws.Cells.Item(83, 14).Value = -29984
$ws.Cells.Item(132, 8).Value = 3075.8076
$ws.Cells.Item(132, 9).Value = 3058.6155
$ws.Cells.Item(132, 11).Value = 9175.8465
$ws.Cells.Item(132, 13).Value = -6645.8465
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4770.7354
$ws.Cells.Item(122, 9).Value = 4627.273
$ws.Cells.Item(122, 10).Value = 5033.75
$ws.Cells.Item(122, 11).Value = 13881.819
$ws.Cells.Item(122, 12).Value = 15101.25
$ws.Cells.Item(122, 13).Value = -11431.819
$ws.Cells.Item(122, 14).Value = -20001.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109, 8).Value = 71400
$ws.Cells.Item(109, 10).Value = 71400
$ws.Cells.Item(109, 12).Value = 71400
$ws.Cells.Item(109, 14).Value = -74174
